$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows where Speaker (column D) is "Davis" -> rename to "T"
$davisRows = @(2,4,10,11,12,13,14,15,16,17,19,21,58,59,61,63,64,65,69,70,71,73,74,76,83,85,86,125,126,127,128,146,147)

# Rows where Speaker (column D) is "Student" -> rename to "S"
$studentRows = @(9,29,30,44,45,48,50,60,67,75,77,79,80,93,135,148,157,160)

foreach ($r in $davisRows) {
    $ws.Range("D$r").Value = "T"
}

foreach ($r in $studentRows) {
    $ws.Range("D$r").Value = "S"
}

# Row 76 also has "Davis" at the start of the Sentence text in column E
$ws.Range("E76").Value = "T I guess, maybe I will."
